$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters with their updated values for rows 2 and 3 (both rows share the same new values)
$updates = @{
    "D"  = -0.0958
    "I"  = -0.04049346879535559
    "J"  = -0.04049346879535559
    "K"  = -0.709
    "L"  = -0.1029027576197388
    "U"  = 0.032
    "V"  = 0.008101265822784811
    "W"  = 0.1729268292682927
    "X"  = 0.2736941782819743
    "Y"  = -0.1007673490136816
    "Z"  = 0.6405131542251556
    "AA" = -0.02593659942363112
    "AB" = 0.06684520700920227
    "AC" = -0.09278180643283339
    "AD" = 21.3
    "AF" = 21.3
    "AG" = 21.268
    "AH" = 0.8435643564356436
    "AI" = 1.287787182587666
    "AJ" = 0.8433658497898326
    "AK" = 1.288345044826751
    "AL" = 0.426
    "AM" = 0.426
    "AO" = -0.6549295774647889
    "AQ" = -0.6549295774647889
}

foreach ($row in 2..3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
